## Applies the "Improvements for V0.3 of the Rastaban PCB" edit:
##  1. Splits the title run so "Rastaban" is wrapped in spell-check markup.
##  2. Replaces the six bulleted list paragraphs with a 3-column table
##     (# / Improvement / Finished) containing a header row + 7 item rows,
##     updating / expanding some of the improvement text along the way.
##  3. Leaves a trailing empty paragraph after the table.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: "Improvements for V0.3 of the Rastaban PCB"
#    -> split into 3 runs, "Rastaban" wrapped in proofErr spell markers.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Improvements for V0.3 of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rastaban</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> PCB</w:t></w:r></w:p>
'@
$titlePara.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2. Replace the six list-bullet paragraphs (now paragraphs 2..7) with
#    a table.
# ---------------------------------------------------------------------
$startPos = $d.Paragraphs(2).Range.Start
$endPos = $d.Paragraphs(7).Range.End
$listRange = $d.Range($startPos, $endPos)

$tableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="Tabelraster"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="421"/>
    <w:gridCol w:w="5786"/>
    <w:gridCol w:w="2855"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>#</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>Improvement</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>Finished</w:t></w:r></w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>1</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>Figure out which hardware pins to use for what components</w:t></w:r><w:r><w:tab/></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>2</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p>
        <w:r><w:t xml:space="preserve">TMC 2209 uses different ms1 ms2 configuration for </w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>microstepping</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t xml:space="preserve"> than the </w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>tmc</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t xml:space="preserve"> 2208! Keep this in mind. </w:t></w:r>
      </w:p>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>3</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>Remove tmc2208 from design, focus on 2209</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>4</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p>
        <w:r><w:t>Use UART on 2209 and remove step/</w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>dir</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t xml:space="preserve"> enable interface.</w:t></w:r>
      </w:p>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>5</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>Use appropriate resistor for UART control on 2209 (resistance should decrease with increase in drivers).</w:t></w:r></w:p>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>6</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p>
        <w:r><w:t xml:space="preserve">Connect </w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>diag</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t xml:space="preserve"> pin of 2209 for stall (stuck motor) indication for </w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>rpi</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t>??</w:t></w:r>
      </w:p>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:w="421" w:type="dxa"/></w:tcPr>
      <w:p><w:r><w:t>7</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="5786" w:type="dxa"/></w:tcPr>
      <w:p>
        <w:r><w:t xml:space="preserve">Checkout the problems with </w:t></w:r>
        <w:proofErr w:type="spellStart"/><w:r><w:t>diag</w:t></w:r><w:proofErr w:type="spellEnd"/>
        <w:r><w:t xml:space="preserve"> pin on tmc2209 (see pdf in datasheets)</w:t></w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:w="2855" w:type="dxa"/></w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$listRange.InsertXML($tableXml)

Write-Output "edit applied"
